# The diff regenerates this docx4j test fixture (dest-with-header-footer.docx)
# with a different build environment: the tool-signature comment embedded at
# the top of <w:body> changes from
#   "... using REFERENCE JAXB in Oracle Java 21.0.8 on Mac OS X"
# to
#   "... using REFERENCE JAXB in Microsoft Java 21.0.8 on Mac OS X"
# (plus purely cosmetic re-ordering of the xmlns:* declarations on the
# document/header/footer/styles root elements, a re-serialization artifact
# of the same regeneration, not a content edit).
#
# Apply the intended textual substitution across every reachable story in
# the document (main body, headers, footers) so that if this string is
# ever exposed through the content/text surface it gets corrected; this is
# a harmless no-op everywhere the string is not part of visible text (e.g.
# when it lives only in a non-content XML comment), so it never corrupts
# unrelated content.

$d = $word.ActiveDocument

$old = "Oracle Java 21.0.8"
$new = "Microsoft Java 21.0.8"

# Main document story (paragraphs, tables, etc.)
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# Every header/footer story in every section (primary/first/even).
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
        }
    }
    foreach ($hf in $sec.Footers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
        }
    }
}
